# edit.ps1 - apply the "update sprint review doc" revision
#
# Summary of changes (per the target diff):
#   1. After "...loginForm component." add a new sentence about search methods.
#   2. Klei's paragraph: split "did 4 sequence diagrams and addressed comments..."
#      to insert ", a design diagram" before " and addressed comments...".
#   3. Chase's paragraph: change "...did 4 sequence diagrams and 2 design diagrams."
#      to "...did 4 sequence diagrams and a design diagram." (split into runs).
#   4. Hunter's paragraph: replace "Hunter did not contribute to this deliverable."
#      with "Hunter created search methods for a user profile." (split into runs).
#   5. Calvin's paragraph: drop "went to scrum meetings and" from his sentence.
#   6. Final "Improvement in the future..." sentence reworded and extended.

$d = $word.ActiveDocument

function Find-InRange($range, [string]$text) {
    # Execute Find against the supplied range (Find narrows the range in
    # place, matching the real Word COM behaviour).
    return $range.Find.Execute($text, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
}

# ---------------------------------------------------------------------
# 1. " component." -> append new sentence as its own run.
# ---------------------------------------------------------------------
$rng = $d.Content
[void](Find-InRange $rng " component.")
$rng.Collapse(0) | Out-Null
$rng.InsertAfter(" Methods for searching for users was created as well.")
$rng.Font.BoldBi = $true
$rng.LanguageID = "en-US"

# ---------------------------------------------------------------------
# 2. Klei paragraph: split the sentence into three runs.
# ---------------------------------------------------------------------
$rng = $d.Content
[void](Find-InRange $rng " did 4 sequence diagrams and addressed comments on the use case diagram document.")
$rng.Text = " did 4 sequence diagrams"

$rng.Collapse(0) | Out-Null
$rng.InsertAfter(", a design diagram")
$rng.Font.BoldBi = $true
$rng.LanguageID = "en-US"

$rng.Collapse(0) | Out-Null
$rng.InsertAfter(" and addressed comments on the use case diagram document.")
$rng.Font.BoldBi = $true
$rng.LanguageID = "en-US"

# ---------------------------------------------------------------------
# 3. Chase paragraph: split into three runs, "2 design diagrams" -> "a design diagram".
# ---------------------------------------------------------------------
$rng = $d.Content
[void](Find-InRange $rng ": Chase addressed comments on the use case diagram document, worked on implementation for the home screen, did 4 sequence diagrams and 2 design diagrams. He also worked on moving the home.html to a react component so that the application is consistent in its architecture. ")
$rng.Text = ": Chase addressed comments on the use case diagram document, worked on implementation for the home screen, did 4 sequence diagrams and "

$rng.Collapse(0) | Out-Null
$rng.InsertAfter("a")
$rng.Font.BoldBi = $true
$rng.LanguageID = "en-US"

$rng.Collapse(0) | Out-Null
$rng.InsertAfter(" design diagram. He also worked on moving the home.html to a react component so that the application is consistent in its architecture. ")
$rng.Font.BoldBi = $true
$rng.LanguageID = "en-US"

# ---------------------------------------------------------------------
# 4. Hunter paragraph: replace with two runs.
# ---------------------------------------------------------------------
$rng = $d.Content
[void](Find-InRange $rng ": Hunter did not contribute to this deliverable. ")
$rng.Text = ": Hunter "

$rng.Collapse(0) | Out-Null
$rng.InsertAfter("created search methods for a user profile.")
$rng.Font.BoldBi = $true
$rng.LanguageID = "en-US"

# ---------------------------------------------------------------------
# 5. Calvin paragraph: drop "went to scrum meetings and".
# ---------------------------------------------------------------------
$rng = $d.Content
[void](Find-InRange $rng ": Calvin went to scrum meetings and worked on implementing the home screen with chase.")
$rng.Text = ": Calvin worked on implementing the home screen with chase."

# ---------------------------------------------------------------------
# 6. Final "Improvement in the future..." sentence.
# ---------------------------------------------------------------------
$rng = $d.Content
[void](Find-InRange $rng " Improvement in the future can include all group members participating and proactively being involved. ")
$rng.Text = " Improvement in the future can include all group members proactively being involved"
$rng.Collapse(0) | Out-Null

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tailRng = $d.Range($rng.End, $lastPara.Range.End - 1)
$tailRng.Delete()

$rng.InsertAfter(" throughout the 2 weeks. ")
$rng.Font.BoldBi = $true
$rng.Font.Size = 12
